# Apply crypto price/volume updates per the Mon May 27 21:55:52 UTC 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "69.622.78") that must stay text, not be
# auto-converted to a number by Excels input parser.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextCell $ws.Range("D2") "69.622.78"
$ws.Range("E2").Value = "  +1.79%  "

Set-TextCell $ws.Range("D3") "3.885.92"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.08%  "

Set-TextCell $ws.Range("D5") "603.95"
$ws.Range("E5").Value = "  +0.82%  "

Set-TextCell $ws.Range("D6") "170.14"
$ws.Range("E6").Value = "  +4.74%  "

Set-TextCell $ws.Range("D7") "3.889.98"
$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("E12").Value = "  +1.83%  "

$ws.Range("E13").Value = "  +4.92%  "

Set-TextCell $ws.Range("D14") "38.15"
$ws.Range("E14").Value = "  +3.81%  "

Set-TextCell $ws.Range("D15") "4.539.92"
$ws.Range("E15").Value = "  +1.33%  "

Set-TextCell $ws.Range("D16") "3.877.89"
$ws.Range("E16").Value = "  +1.21%  "

Set-TextCell $ws.Range("D17") "69.659.13"
$ws.Range("E17").Value = "  +1.53%  "

Set-TextCell $ws.Range("D18") "18.68"
$ws.Range("E18").Value = "  +9.33%  "

Set-TextCell $ws.Range("D19") "7.64"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("E20").Value = "  -0.80%  "

Set-TextCell $ws.Range("D21") "11.15"
$ws.Range("E21").Value = "  -0.88%  "

Set-TextCell $ws.Range("D22") "489.18"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("E23").Value = "  +4.39%  "

Set-TextCell $ws.Range("D24") "0.0000164"
$ws.Range("E24").Value = "  +2.41%  "

Set-TextCell $ws.Range("D25") "85.21"
$ws.Range("E25").Value = "  +1.63%  "

$ws.Range("E26").Value = "  +3.83%  "

Set-TextCell $ws.Range("D27") "12.36"
$ws.Range("E27").Value = "  +2.34%  "

Set-TextCell $ws.Range("D28") "10.12"
$ws.Range("E28").Value = "  +2.11%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  +1.13%  "

Set-TextCell $ws.Range("D31") "4.035.42"
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("E32").Value = "  +2.17%  "

Set-TextCell $ws.Range("D33") "7.86"
$ws.Range("E33").Value = "  +0.53%  "

Set-TextCell $ws.Range("D34") "31.83"
$ws.Range("E34").Value = "  -0.46%  "

Set-TextCell $ws.Range("D35") "3.851.89"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("E36").Value = "  +0.51%  "

$ws.Range("E37").Value = "  +4.45%  "

$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D38") "1.03"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D39") "0.142"
$ws.Range("E39").Value = "  +2.14%  "

Set-TextCell $ws.Range("D40") "3.37"
$ws.Range("E40").Value = "  +14.17%  "

Set-TextCell $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("E42").Value = "  +3.41%  "

Set-TextCell $ws.Range("D43") "2.07"
$ws.Range("E43").Value = "  +5.54%  "

Set-TextCell $ws.Range("D44") "437.12"
$ws.Range("E44").Value = "  +1.73%  "

Set-TextCell $ws.Range("D45") "48.11"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("E46").Value = "  +3.96%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("E48").Value = "  +22.99%  "

Set-TextCell $ws.Range("D49") "0.0366"
$ws.Range("E49").Value = "  +3.24%  "

Set-TextCell $ws.Range("D50") "144.06"
$ws.Range("E50").Value = "  +0.77%  "

Set-TextCell $ws.Range("D51") "40.53"
$ws.Range("E51").Value = "  +4.79%  "
